$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

$tr.Characters(1, 84).Text = '[MARKER_OPENING 0 -ERR:REF-NOT-FOUND-]   =andray()[MARKER_CLOSING 1 &lt;/a:t>&lt;/a:r>]   '

# Shape 2
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$tr.Characters(8118, 87).Text = '[MARKER_OPENING 190 -ERR:REF-NOT-FOUND-]    ogday. [MARKER_CLOSING 191 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(8033, 85).Text = '[MARKER_OPENING 188 -ERR:REF-NOT-FOUND-]   azylay[MARKER_CLOSING 189 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(7939, 94).Text = '[MARKER_OPENING 186 -ERR:REF-NOT-FOUND-]    overhay ethay [MARKER_CLOSING 187 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(7853, 86).Text = '[MARKER_OPENING 184 -ERR:REF-NOT-FOUND-]   umpsjay[MARKER_CLOSING 185 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(7767, 86).Text = '[MARKER_OPENING 182 -ERR:REF-NOT-FOUND-]    oxfay [MARKER_CLOSING 183 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(7681, 86).Text = '[MARKER_OPENING 180 -ERR:REF-NOT-FOUND-]   ownbray[MARKER_CLOSING 181 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(7580, 101).Text = '[MARKER_OPENING 178 -ERR:REF-NOT-FOUND-]    ogday. ethay uickqay [MARKER_CLOSING 179 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(7495, 85).Text = '[MARKER_OPENING 176 -ERR:REF-NOT-FOUND-]   azylay[MARKER_CLOSING 177 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(7401, 94).Text = '[MARKER_OPENING 174 -ERR:REF-NOT-FOUND-]    overhay ethay [MARKER_CLOSING 175 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(7315, 86).Text = '[MARKER_OPENING 172 -ERR:REF-NOT-FOUND-]   umpsjay[MARKER_CLOSING 173 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(7229, 86).Text = '[MARKER_OPENING 170 -ERR:REF-NOT-FOUND-]    oxfay [MARKER_CLOSING 171 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(7143, 86).Text = '[MARKER_OPENING 168 -ERR:REF-NOT-FOUND-]   ownbray[MARKER_CLOSING 169 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(7042, 101).Text = '[MARKER_OPENING 166 -ERR:REF-NOT-FOUND-]    ogday. ethay uickqay [MARKER_CLOSING 167 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(6957, 85).Text = '[MARKER_OPENING 164 -ERR:REF-NOT-FOUND-]   azylay[MARKER_CLOSING 165 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(6863, 94).Text = '[MARKER_OPENING 162 -ERR:REF-NOT-FOUND-]    overhay ethay [MARKER_CLOSING 163 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(6777, 86).Text = '[MARKER_OPENING 160 -ERR:REF-NOT-FOUND-]   umpsjay[MARKER_CLOSING 161 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(6691, 86).Text = '[MARKER_OPENING 158 -ERR:REF-NOT-FOUND-]    oxfay [MARKER_CLOSING 159 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(6605, 86).Text = '[MARKER_OPENING 156 -ERR:REF-NOT-FOUND-]   ownbray[MARKER_CLOSING 157 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(6504, 101).Text = '[MARKER_OPENING 154 -ERR:REF-NOT-FOUND-]    ogday. ethay uickqay [MARKER_CLOSING 155 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(6419, 85).Text = '[MARKER_OPENING 152 -ERR:REF-NOT-FOUND-]   azylay[MARKER_CLOSING 153 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(6325, 94).Text = '[MARKER_OPENING 150 -ERR:REF-NOT-FOUND-]    overhay ethay [MARKER_CLOSING 151 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(6239, 86).Text = '[MARKER_OPENING 148 -ERR:REF-NOT-FOUND-]   umpsjay[MARKER_CLOSING 149 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(6153, 86).Text = '[MARKER_OPENING 146 -ERR:REF-NOT-FOUND-]    oxfay [MARKER_CLOSING 147 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(6067, 86).Text = '[MARKER_OPENING 144 -ERR:REF-NOT-FOUND-]   ownbray[MARKER_CLOSING 145 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(5966, 101).Text = '[MARKER_OPENING 142 -ERR:REF-NOT-FOUND-]    ogday. ethay uickqay [MARKER_CLOSING 143 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(5881, 85).Text = '[MARKER_OPENING 140 -ERR:REF-NOT-FOUND-]   azylay[MARKER_CLOSING 141 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(5787, 94).Text = '[MARKER_OPENING 138 -ERR:REF-NOT-FOUND-]    overhay ethay [MARKER_CLOSING 139 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(5701, 86).Text = '[MARKER_OPENING 136 -ERR:REF-NOT-FOUND-]   umpsjay[MARKER_CLOSING 137 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(5615, 86).Text = '[MARKER_OPENING 134 -ERR:REF-NOT-FOUND-]    oxfay [MARKER_CLOSING 135 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(5529, 86).Text = '[MARKER_OPENING 132 -ERR:REF-NOT-FOUND-]   ownbray[MARKER_CLOSING 133 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(5436, 93).Text = '[MARKER_OPENING 130 -ERR:REF-NOT-FOUND-]   ethay uickqay [MARKER_CLOSING 131 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(5347, 87).Text = '[MARKER_OPENING 125 -ERR:REF-NOT-FOUND-]    ogday. [MARKER_CLOSING 126 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(5262, 85).Text = '[MARKER_OPENING 123 -ERR:REF-NOT-FOUND-]   azylay[MARKER_CLOSING 124 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(5168, 94).Text = '[MARKER_OPENING 121 -ERR:REF-NOT-FOUND-]    overhay ethay [MARKER_CLOSING 122 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(5082, 86).Text = '[MARKER_OPENING 119 -ERR:REF-NOT-FOUND-]   umpsjay[MARKER_CLOSING 120 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(4996, 86).Text = '[MARKER_OPENING 117 -ERR:REF-NOT-FOUND-]    oxfay [MARKER_CLOSING 118 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(4910, 86).Text = '[MARKER_OPENING 115 -ERR:REF-NOT-FOUND-]   ownbray[MARKER_CLOSING 116 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(4809, 101).Text = '[MARKER_OPENING 113 -ERR:REF-NOT-FOUND-]    ogday. ethay uickqay [MARKER_CLOSING 114 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(4724, 85).Text = '[MARKER_OPENING 111 -ERR:REF-NOT-FOUND-]   azylay[MARKER_CLOSING 112 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(4630, 94).Text = '[MARKER_OPENING 109 -ERR:REF-NOT-FOUND-]    overhay ethay [MARKER_CLOSING 110 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(4544, 86).Text = '[MARKER_OPENING 107 -ERR:REF-NOT-FOUND-]   umpsjay[MARKER_CLOSING 108 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(4458, 86).Text = '[MARKER_OPENING 105 -ERR:REF-NOT-FOUND-]    oxfay [MARKER_CLOSING 106 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(4372, 86).Text = '[MARKER_OPENING 103 -ERR:REF-NOT-FOUND-]   ownbray[MARKER_CLOSING 104 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(4271, 101).Text = '[MARKER_OPENING 101 -ERR:REF-NOT-FOUND-]    ogday. ethay uickqay [MARKER_CLOSING 102 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(4187, 84).Text = '[MARKER_OPENING 99 -ERR:REF-NOT-FOUND-]   azylay[MARKER_CLOSING 100 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(4095, 92).Text = '[MARKER_OPENING 97 -ERR:REF-NOT-FOUND-]    overhay ethay [MARKER_CLOSING 98 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(4011, 84).Text = '[MARKER_OPENING 95 -ERR:REF-NOT-FOUND-]   umpsjay[MARKER_CLOSING 96 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(3927, 84).Text = '[MARKER_OPENING 93 -ERR:REF-NOT-FOUND-]    oxfay [MARKER_CLOSING 94 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(3843, 84).Text = '[MARKER_OPENING 91 -ERR:REF-NOT-FOUND-]   ownbray[MARKER_CLOSING 92 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(3744, 99).Text = '[MARKER_OPENING 89 -ERR:REF-NOT-FOUND-]    ogday. ethay uickqay [MARKER_CLOSING 90 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(3661, 83).Text = '[MARKER_OPENING 87 -ERR:REF-NOT-FOUND-]   azylay[MARKER_CLOSING 88 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(3569, 92).Text = '[MARKER_OPENING 85 -ERR:REF-NOT-FOUND-]    overhay ethay [MARKER_CLOSING 86 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(3485, 84).Text = '[MARKER_OPENING 83 -ERR:REF-NOT-FOUND-]   umpsjay[MARKER_CLOSING 84 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(3401, 84).Text = '[MARKER_OPENING 81 -ERR:REF-NOT-FOUND-]    oxfay [MARKER_CLOSING 82 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(3317, 84).Text = '[MARKER_OPENING 79 -ERR:REF-NOT-FOUND-]   ownbray[MARKER_CLOSING 80 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(3218, 99).Text = '[MARKER_OPENING 77 -ERR:REF-NOT-FOUND-]    ogday. ethay uickqay [MARKER_CLOSING 78 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(3135, 83).Text = '[MARKER_OPENING 75 -ERR:REF-NOT-FOUND-]   azylay[MARKER_CLOSING 76 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(3043, 92).Text = '[MARKER_OPENING 73 -ERR:REF-NOT-FOUND-]    overhay ethay [MARKER_CLOSING 74 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(2959, 84).Text = '[MARKER_OPENING 71 -ERR:REF-NOT-FOUND-]   umpsjay[MARKER_CLOSING 72 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(2875, 84).Text = '[MARKER_OPENING 69 -ERR:REF-NOT-FOUND-]    oxfay [MARKER_CLOSING 70 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(2791, 84).Text = '[MARKER_OPENING 67 -ERR:REF-NOT-FOUND-]   ownbray[MARKER_CLOSING 68 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(2700, 91).Text = '[MARKER_OPENING 65 -ERR:REF-NOT-FOUND-]   ethay uickqay [MARKER_CLOSING 66 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(2613, 85).Text = '[MARKER_OPENING 60 -ERR:REF-NOT-FOUND-]    ogday. [MARKER_CLOSING 61 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(2530, 83).Text = '[MARKER_OPENING 58 -ERR:REF-NOT-FOUND-]   azylay[MARKER_CLOSING 59 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(2438, 92).Text = '[MARKER_OPENING 56 -ERR:REF-NOT-FOUND-]    overhay ethay [MARKER_CLOSING 57 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(2354, 84).Text = '[MARKER_OPENING 54 -ERR:REF-NOT-FOUND-]   umpsjay[MARKER_CLOSING 55 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(2270, 84).Text = '[MARKER_OPENING 52 -ERR:REF-NOT-FOUND-]    oxfay [MARKER_CLOSING 53 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(2186, 84).Text = '[MARKER_OPENING 50 -ERR:REF-NOT-FOUND-]   ownbray[MARKER_CLOSING 51 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(2087, 99).Text = '[MARKER_OPENING 48 -ERR:REF-NOT-FOUND-]    ogday. ethay uickqay [MARKER_CLOSING 49 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(2004, 83).Text = '[MARKER_OPENING 46 -ERR:REF-NOT-FOUND-]   azylay[MARKER_CLOSING 47 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(1912, 92).Text = '[MARKER_OPENING 44 -ERR:REF-NOT-FOUND-]    overhay ethay [MARKER_CLOSING 45 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(1828, 84).Text = '[MARKER_OPENING 42 -ERR:REF-NOT-FOUND-]   umpsjay[MARKER_CLOSING 43 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(1744, 84).Text = '[MARKER_OPENING 40 -ERR:REF-NOT-FOUND-]    oxfay [MARKER_CLOSING 41 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(1660, 84).Text = '[MARKER_OPENING 38 -ERR:REF-NOT-FOUND-]   ownbray[MARKER_CLOSING 39 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(1561, 99).Text = '[MARKER_OPENING 36 -ERR:REF-NOT-FOUND-]    ogday. ethay uickqay [MARKER_CLOSING 37 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(1478, 83).Text = '[MARKER_OPENING 34 -ERR:REF-NOT-FOUND-]   azylay[MARKER_CLOSING 35 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(1386, 92).Text = '[MARKER_OPENING 32 -ERR:REF-NOT-FOUND-]    overhay ethay [MARKER_CLOSING 33 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(1302, 84).Text = '[MARKER_OPENING 30 -ERR:REF-NOT-FOUND-]   umpsjay[MARKER_CLOSING 31 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(1218, 84).Text = '[MARKER_OPENING 28 -ERR:REF-NOT-FOUND-]    oxfay [MARKER_CLOSING 29 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(1134, 84).Text = '[MARKER_OPENING 26 -ERR:REF-NOT-FOUND-]   ownbray[MARKER_CLOSING 27 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(1035, 99).Text = '[MARKER_OPENING 24 -ERR:REF-NOT-FOUND-]    ogday. ethay uickqay [MARKER_CLOSING 25 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(952, 83).Text = '[MARKER_OPENING 22 -ERR:REF-NOT-FOUND-]   azylay[MARKER_CLOSING 23 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(860, 92).Text = '[MARKER_OPENING 20 -ERR:REF-NOT-FOUND-]    overhay ethay [MARKER_CLOSING 21 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(776, 84).Text = '[MARKER_OPENING 18 -ERR:REF-NOT-FOUND-]   umpsjay[MARKER_CLOSING 19 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(692, 84).Text = '[MARKER_OPENING 16 -ERR:REF-NOT-FOUND-]    oxfay [MARKER_CLOSING 17 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(608, 84).Text = '[MARKER_OPENING 14 -ERR:REF-NOT-FOUND-]   ownbray[MARKER_CLOSING 15 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(509, 99).Text = '[MARKER_OPENING 12 -ERR:REF-NOT-FOUND-]    ogday. ethay uickqay [MARKER_CLOSING 13 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(426, 83).Text = '[MARKER_OPENING 10 -ERR:REF-NOT-FOUND-]   azylay[MARKER_CLOSING 11 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(336, 90).Text = '[MARKER_OPENING 8 -ERR:REF-NOT-FOUND-]    overhay ethay [MARKER_CLOSING 9 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(254, 82).Text = '[MARKER_OPENING 6 -ERR:REF-NOT-FOUND-]   umpsjay[MARKER_CLOSING 7 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(172, 82).Text = '[MARKER_OPENING 4 -ERR:REF-NOT-FOUND-]    oxfay [MARKER_CLOSING 5 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(90, 82).Text = '[MARKER_OPENING 2 -ERR:REF-NOT-FOUND-]   ownbray[MARKER_CLOSING 3 &lt;/a:t>&lt;/a:r>]   '
$tr.Characters(1, 89).Text = '[MARKER_OPENING 0 -ERR:REF-NOT-FOUND-]   ethay uickqay [MARKER_CLOSING 1 &lt;/a:t>&lt;/a:r>]   '
